$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8068333864212036
$ws.Range("B1").Value = 3.668709993362427
$ws.Range("C1").Value = 3.412396907806396
$ws.Range("D1").Value = 2.853431463241577
$ws.Range("E1").Value = 1.766730308532715
